# Updates the "想去人数" (want-to-go count) column F values across the
# four worksheets, matching the re-scraped numbers from the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitValues = @{
    3  = 152
    4  = 1783
    5  = 3333
    6  = 1088
    7  = 2223
    8  = 2133
    10 = 609
    12 = 1682
    16 = 99
    17 = 221
    18 = 1594
    19 = 641
    20 = 732
    21 = 616
    22 = 12288
    23 = 12347
    24 = 913
    25 = 704
    27 = 42
    28 = 28
    29 = 380
    30 = 1926
    31 = 2
    32 = 5
    33 = 202
    34 = 599
}
foreach ($row in $exhibitValues.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitValues[$row]
}

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$showValues = @{
    6 = 42
    7 = 37
}
foreach ($row in $showValues.Keys) {
    $wsShow.Cells.Item($row, 6).Value = $showValues[$row]
}

# Sheet "本地生活" (Local life)
$wsLocal = $wb.Worksheets.Item("本地生活")
$localValues = @{
    2 = 77
}
foreach ($row in $localValues.Keys) {
    $wsLocal.Cells.Item($row, 6).Value = $localValues[$row]
}

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$allValues = @{
    3  = 77
    4  = 152
    5  = 1783
    6  = 3333
    7  = 1088
    8  = 2223
    9  = 2133
    11 = 609
    13 = 1682
    19 = 99
    21 = 221
    22 = 1594
    23 = 641
    24 = 732
    25 = 616
    26 = 12288
    27 = 12347
    28 = 913
    29 = 704
    31 = 42
    32 = 28
    33 = 380
    34 = 1926
    35 = 2
    37 = 5
    38 = 42
    39 = 202
    40 = 599
    41 = 37
}
foreach ($row in $allValues.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allValues[$row]
}
